# GuildConfig.xlsx: unify the conception of DataNode, DataTable, Entity.
# Rename the sheet from the old "Property1" label to "DataNode", and carry
# over the row-height/selection changes that came along with the resave.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet: Property1 -> DataNode
$ws.Name = "DataNode"

# 2. Row-height adjustments (sheet default row height moved from 15 to 13.5,
#    the header row grew to fit wrapped text, and the sub-header row shrank).
$ws.Rows.Item(1).RowHeight = 40.5   # title/header row
$ws.Rows.Item(2).RowHeight = 13.5
$ws.Rows.Item(3).RowHeight = 13.5
$ws.Rows.Item(4).RowHeight = 13.5
$ws.Rows.Item(5).RowHeight = 13.5
$ws.Rows.Item(6).RowHeight = 13.5
$ws.Rows.Item(7).RowHeight = 13.5
$ws.Rows.Item(8).RowHeight = 27     # column-description row
$ws.Rows.Item(9).RowHeight = 13.5
$ws.Rows.Item(10).RowHeight = 13.5
$ws.Rows.Item(11).RowHeight = 13.5

# 3. Selection moved off the frozen pane's corner cell to D36.
$ws.Range("D36").Select() | Out-Null
